# Cập nhật kết nối online lên heroku
# Fill in the "online" (column C) connection settings on the "Kết nối CSDL"
# sheet with the new Heroku / ClearDB values, then leave that sheet active
# (matching where the user was last working).

$wb = $excel.ActiveWorkbook
$dbSheet = $wb.Worksheets.Item("Kết nối CSDL")

# user
$dbSheet.Range("C4").Value = "b20a7ec83541b0"

# password
$dbSheet.Range("C5").Value = "64ffe227"

# host
$dbSheet.Range("C2").Value = "us-cdbr-iron-east-03.cleardb.net"

# port - copy B3's number format/style onto C3, then set the value
$dbSheet.Range("B3").Copy()
$dbSheet.Range("C3").PasteSpecial(-4122)
$dbSheet.Range("C3").Value = 3306

# database
$dbSheet.Range("C6").Value = "heroku_633ee9287d27a78"

# The user ended their session on the "Kết nối CSDL" sheet, with D16 selected
$dbSheet.Activate()
$dbSheet.Range("D16").Select()
